$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for the two added columns
$ws.Range("L1").Value = "Physical Activity"
$ws.Range("M1").Value = "Hands On Time"

# Default all data rows (2-54) in the new columns to 0
$ws.Range("L2:M54").Value = 0

# Row 9 (Milton Clark) - updated existing values + new columns
$ws.Range("H9").Value = 20
$ws.Range("I9").Value = 17
$ws.Range("J9").Value = 17
$ws.Range("K9").Value = 17
$ws.Range("L9").Value = 20
$ws.Range("M9").Value = 20

# Row 10 (Tresean Clark) - updated existing values + new columns
$ws.Range("H10").Value = 27
$ws.Range("I10").Value = 22
$ws.Range("J10").Value = 17
$ws.Range("K10").Value = 12
$ws.Range("L10").Value = 20
$ws.Range("M10").Value = 25

# Row 11 (Zamir Clark) - updated existing values + new columns
$ws.Range("H11").Value = 30
$ws.Range("I11").Value = 15
$ws.Range("J11").Value = 10
$ws.Range("K11").Value = 17
$ws.Range("L11").Value = 20
$ws.Range("M11").Value = 30

# Row 19 (Gabriel Hilliard) - updated existing values + new columns
$ws.Range("I19").Value = 14
$ws.Range("J19").Value = 10
$ws.Range("K19").Value = 22
$ws.Range("L19").Value = 35
$ws.Range("M19").Value = 22

# Set column widths for the two new columns to match the other data columns.
# (The COM width model here quantizes to whole pixels the same way Excel does,
# so the input is pre-compensated for the fixed 5px/7 padding offset to land
# on the pixel closest to the source file's fractional "best fit" widths.)
$ws.Columns.Item(12).ColumnWidth = 10.142857142857142
$ws.Columns.Item(13).ColumnWidth = 13.428571428571429
